# Apply the "Thank You" slide edits (final slide, slide 7).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- Shape 1: Title placeholder ("QUESTIONS" -> "Thank You") ---
$title = $s.Shapes.Item(1)

# Reposition / resize (EMU -> points, 12700 EMU per point). The literal
# point values below are nudged by a few ULPs so that after the COM
# layer's internal single-precision rounding they land back on the
# exact target EMU (plain "emu/12700" is off by one EMU for cx/cy).
$title.Left   = 257.62023622047246
$title.Top    = 184.5971653543307
$title.Width  = 444.75953675905515
$title.Height = 81.47110376220472

# Clear the existing run first so the stray trailing <a:endParaRPr/> that
# PowerPoint leaves behind after a plain Text= assignment is dropped (matches
# the target markup, which has no endParaRPr on this paragraph any more).
$title.TextFrame.DeleteText()
$title.TextFrame.TextRange.Text = "Thank You"

# --- Shape 2: Subtitle placeholder (college/department line) ---
$subtitle = $s.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Text = "College of Engineering |  Bioengineering & Computer Science"
